$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the closing "italic description" paragraph's text (was the
#    meta-description blurb, becomes the AI image-prompt text). Do this
#    first (before any earlier-document edits) so paragraph indices for
#    the tail of the document stay stable while we touch them.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Text = "Prompt: Please create a feature image for the slot game ""Double Happiness"" that fits the game's theme and features a happy Maya warrior with glasses. The image should be in cartoon style."

# ---------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Double Happiness Slot for Free -
#    Review 2021" paragraph that used to sit right before it (its
#    content has moved up to become the new "Meta description"
#    paragraph near the top of the document).
# ---------------------------------------------------------------------
$dupIndex = $count - 1
$dupPara = $d.Paragraphs($dupIndex)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Insert the new "Meta description" paragraph right after the title
#    (first) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$fullText = "Meta description: Experience the Chinese-themed Double Happiness slot for free. Read our review on the game's features, symbols, and winning probability."
$ins = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$ins.InsertAfter($fullText)

$labelStart = $d.Paragraphs(2).Range.Start
$labelLen = "Meta description".Length
$boldRange = $d.Range($labelStart, $labelStart + $labelLen)
$boldRange.Bold = 1
